$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "尿微量蛋白"
$ws.Range("C2").Value = "mg/L"
$ws.Range("E2").Value = "0-20"

$ws.Range("A3").Value = "尿糖"
$ws.Range("E3").Value = "-"

$ws.Range("A4").Value = "胆红素"
$ws.Range("E4").Value = "-"

$ws.Range("A5").Value = "尿萌体"
$ws.Range("E5").Value = "-"

$ws.Range("A6").Value = "尿蛋白"
$ws.Range("E6").Value = "-"

$ws.Range("A7").Value = "亚硝酸盐"
$ws.Range("E7").Value = "-"

$ws.Range("A8").Value = "潜血"
$ws.Range("E8").Value = "-"

$ws.Range("A9").Value = "尿白细胞"
$ws.Range("E9").Value = "-"

$ws.Range("A10").Value = "尿胆原"
$ws.Range("E10").Value = "0-2"

$ws.Range("A11").Value = "尿PH"

$ws.Range("A12").Value = "尿比重"
$ws.Range("E12").Value = "1.005-1.030"

$ws.Range("A13").Value = "红细胞"

$ws.Range("A14").Value = "白细胞"
$ws.Range("E14").Value = "0-28"

$ws.Range("A15").Value = "细胞管型"
$ws.Range("E15").Value = "0-3"

$ws.Range("A16").Value = "上皮细胞"
$ws.Range("E16").Value = "0-25"

$ws.Range("A17").Value = "其他结晶"

$ws.Range("A18").Value = "粘液"
$ws.Range("E18").Value = "0-5"
